$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 IP value changes from 127.0.0.1 to 192.168.1.113
$ws.Range("C2").Value = "192.168.1.113"
$ws.Range("C2").NumberFormat = "@"

# Column C gets its own width (15 characters), no longer shares width with column B
$ws.Columns("C").ColumnWidth = 14.285714285714286

# Selection moves to C2
$ws.Range("C2").Select()
